$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 12: end time updated from 0.59375 to 0.625 (14:15 -> 15:00)
$ws.Range("E12").Value = 0.625

# Row 13: new time-registration entry (same role as row 12)
$ws.Range("A13").Value = $ws.Range("A12").Value2
$ws.Range("B13").Value = $ws.Range("B12").Value2
$ws.Range("C13").Value = 43893
$ws.Range("D13").Value = 0.35416666666666669
$ws.Range("E13").Value = 0.60416666666666663

# Row 14: another new time-registration entry (same role)
$ws.Range("A14").Value = $ws.Range("A12").Value2
$ws.Range("B14").Value = $ws.Range("B12").Value2
$ws.Range("C14").Value = 43893
$ws.Range("D14").Value = 0.60416666666666663
$ws.Range("E14").Value = 0.66666666666666663

# Update the active selection to E15
$ws.Range("E15").Select()
